$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.526.02'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.989.79'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '382.16'
$ws.Range('E5').Value = '  +1.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.78'
$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.548'
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.84'
$ws.Range('E10').Value = '  +1.37%  '
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.459.50'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.46'
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.82'
$ws.Range('E15').Value = '  +3.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.003.03'
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.16'
$ws.Range('E17').Value = '  +4.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.999'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '51.513.98'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.09'
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.63'
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.59'
$ws.Range('E23').Value = '  +2.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.86'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('E25').Value = '  +2.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.87'
$ws.Range('E26').Value = '  -3.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.47'
$ws.Range('E27').Value = '  -2.31%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.167'
$ws.Range('E29').Value = '  +1.84%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '26.07'
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.34'
$ws.Range('E32').Value = '  +3.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.70'
$ws.Range('E33').Value = '  +4.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '51.59'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.04'
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.26'
$ws.Range('E38').Value = '  +3.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.83'
$ws.Range('E39').Value = '  +3.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.117'
$ws.Range('E40').Value = '  +1.61%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.84'
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.56'
$ws.Range('E42').Value = '  +3.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '124.81'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('E44').Value = '  +9.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.44'
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.39'
$ws.Range('E46').Value = '  +3.64%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.03'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.043.81'
$ws.Range('E49').Value = '  +2.66%  '
$ws.Range('B50').Value = 'WOONetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.539'
$ws.Range('E50').Value = '  +16.59%  '
$ws.Range('B51').Value = 'BEAM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0332'
$ws.Range('E51').Value = '  +2.55%  '
